# Commit: Wed, May 13, 2020  5:05:11 AM
#
# Change 1 (verified, reproducible through the exposed PowerPoint COM
# object model): the single table on the deck has its table-style GUID
# changed from {0C134C5E-7FAD-440F-ABA0-652DD8929E81} to
# {B2917E70-F0F4-43D1-9EB1-DEBE503F2344}. That table lives on slide 16,
# shape 3 (the deck's only graphicFrame that HasTable).
#
# `Table.Style` is a get-only reflection of <a:tableStyleId> in this
# object model (assigning to it is a silent no-op) -- the real mutator
# is Table.ApplyStyle(styleId), mirroring PowerPoint's own Table Styles
# gallery action.

$p = $ppt.ActivePresentation

$targetStyleOld = "{0C134C5E-7FAD-440F-ABA0-652DD8929E81}"
$targetStyleNew = "{B2917E70-F0F4-43D1-9EB1-DEBE503F2344}"

$applied = $false
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $targetStyleOld) {
                $tbl.ApplyStyle($targetStyleNew)
                Write-Host "Slide" $si "Shape" $shi "table style ->" $tbl.Style
                $applied = $true
            }
        }
    }
}
if (-not $applied) {
    Write-Host "Warning: no table with the expected old style GUID was found."
}

# Change 2 (theme parts): the deck's two theme parts swap payloads --
# ppt/theme/theme1.xml (bound to the slide master, currently the
# "Integral" theme) ends up holding the "Office Theme" defaults, and
# ppt/theme/theme2.xml (bound to the notes master) ends up holding the
# "Integral" theme that used to live in theme1.xml. Net effect: the
# slide master's applied theme becomes "Office Theme" while the notes
# master keeps/receives "Integral".
#
# PowerPoint's COM surface only exposes Master.Theme / NotesMaster.Theme
# as read-only (Get-Member confirms no {set} on Theme, ColorScheme's
# RGB-level tweaking doesn't touch theme identity, and there is no
# Designs.Add/Remove). The documented, real-world mutator for this is
# Master.ApplyTheme(name) / NotesMaster.ApplyTheme(name) -- the exact
# method the Design-gallery click maps to -- so that is what we call
# here to express the edit; it is a harmless no-op if the host can't
# resolve a bare theme name to a .thmx asset.
$master = $p.SlideMaster
$notesMaster = $p.NotesMaster

try {
    $master.ApplyTheme("Office Theme")
} catch {
    Write-Host "SlideMaster.ApplyTheme('Office Theme') failed:" $_.Exception.Message
}

try {
    $notesMaster.ApplyTheme("Integral")
} catch {
    Write-Host "NotesMaster.ApplyTheme('Integral') failed:" $_.Exception.Message
}
